$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Extend the "Tableau1" table by one row (E5:M49 -> E5:M50)
# ---------------------------------------------------------------
$tbl = $ws.ListObjects.Item(1)
$newRow = $tbl.ListRows.Add()

# ---------------------------------------------------------------
# 2. Copy the formatting of the last existing data row (49) onto
#    the freshly added row (50) so styles/number formats match.
# ---------------------------------------------------------------
$ws.Range("E49:M49").Copy() | Out-Null
$ws.Range("E50:M50").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# The "Heure fin" cell is left blank on the new row and, unlike the
# other rows, is not time-formatted (General style) - copy that
# particular look from a cell that already uses it.
$ws.Range("I49").Copy() | Out-Null
$ws.Range("G50").PasteSpecial(-4122) | Out-Null        # xlPasteFormats

$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# 3. Fill in the new row's data
# ---------------------------------------------------------------
$ws.Range("E50").Value = 44287
$ws.Range("F50").Value = 0.34722222222222227
# "Heure fin" (G50) stays empty

$ws.Range("H50").Formula = '=IF(ISBLANK(Tableau1[[#This Row],[Heure Début]]),"",Tableau1[[#This Row],[Heure fin]]-Tableau1[[#This Row],[Heure Début]])'

$ws.Range("I50").Value = "Développement"
$ws.Range("J50").Value = "Implémenter une fonction de log"
$ws.Range("K50").Value = "CPNV"
$ws.Range("L50").Value = "Créé un fichier qui enregistre tout les évènement important ainsi que la date et l'heure de l'évènement"
$ws.Range("M50").Value = "https://www.studytonight.com/c/programs/misc/display-current-date-and-time`nhttps://stackoverflow.com/questions/1442116/how-to-get-the-date-and-time-values-in-a-c-program#:~:text=You%20can%20get%20both%20the,time%20and%20date%20in%20UTC."

# ---------------------------------------------------------------
# 4. Widen the "Source" column (M) so the new long links fit
# ---------------------------------------------------------------
$ws.Columns("M").ColumnWidth = 25.5

# ---------------------------------------------------------------
# 5. Recompute the (word-wrapped) row heights impacted by the
#    wider column / new content
# ---------------------------------------------------------------
$ws.Rows(42).RowHeight = 100.8
$ws.Rows(48).RowHeight = 57.6
$ws.Rows(50).RowHeight = 144

# ---------------------------------------------------------------
# 6. Update the view: scroll down and select the new "Heure fin"
#    cell, like a user having just finished typing the new row.
# ---------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 49
$ws.Range("G50").Select() | Out-Null
